# LU-BU_月度数据.xlsx (Sheet1) — update forecast ("预测值", column C) and one
# actual ("真实值", column B) figure to the refreshed values from the
# coworker notebook re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 461.4    # 2025/12/31 预测值: -71.2 -> 461.4
$ws.Range("C3").Value  = 416.6    # 2025/11/30 预测值: -52.5 -> 416.6
$ws.Range("B4").Value  = -49.9    # 2025/10/31 真实值: -51.2 -> -49.9
$ws.Range("C4").Value  = 108.2    # 2025/10/31 预测值: 29.4  -> 108.2
$ws.Range("C5").Value  = 48.1     # 2025/09/30 预测值: 68.1  -> 48.1
$ws.Range("C6").Value  = 16.6     # 2025/08/31 预测值: 16.5  -> 16.6
$ws.Range("C11").Value = 56.7     # 2025/03/31 预测值: 56.6  -> 56.7
$ws.Range("C12").Value = 150.9    # 2025/02/28 预测值: 151   -> 150.9
$ws.Range("C15").Value = 594.5    # 2024/11/30 预测值: 594.6 -> 594.5
$ws.Range("C16").Value = 773.8    # 2024/10/31 预测值: 773.7 -> 773.8
$ws.Range("C25").Value = 459.7    # 2024/01/31 预测值: 459.8 -> 459.7
